$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (M) mirroring the style pattern of the existing
# year columns (D:L), with values taken from the commit diff.

# Header cell M4: year 2021, bold font (fontId 2) + bottom border (borderId 1)
$ws.Range("M4").Value = 2021
$ws.Range("M4").Font.Bold = $true
$ws.Range("M4").Borders.Item(9).LineStyle = 1
$ws.Range("M4").Borders.Item(9).Weight = -4138

# Data cells M5, M6: regular font (fontId 3), no border
$ws.Range("M5").Value = 93.5
$ws.Range("M6").Value = 96.6

# Data cell M7: regular font (fontId 3), bottom border (borderId 1)
$ws.Range("M7").Value = 98.1
$ws.Range("M7").Borders.Item(9).LineStyle = 1
$ws.Range("M7").Borders.Item(9).Weight = -4138

# Empty cell M3, style matches neighboring thin-border row (borderId 1) used
# throughout row 3.
$ws.Range("M3").Borders.Item(9).LineStyle = 1
$ws.Range("M3").Borders.Item(9).Weight = -4138

# Selection moves as Excel would after entering data down to M7 and tabbing over
$ws.Range("N11").Select()
